$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column D
$ws.Range("D1").Value = "ID2"

# Fill D2:D16 with sequential values starting at 34
$values = 34..48
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# Update selection to match the new active range
$ws.Range("D2:D16").Select()
